$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 160; this shifts the existing rows 160-178 down to 161-179
# and extends the used range to R179 (matching the new <dimension ref="A1:R179"/>).
$ws.Rows.Item(160).Insert()

# Populate the newly inserted row 160 with the new weekly price record.
$ws.Cells.Item(160, 1).Value = 9
$ws.Cells.Item(160, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(160, 3).Value = "Metropolitana"
$ws.Cells.Item(160, 4).Value = 45142
$ws.Cells.Item(160, 5).Value = 13
$ws.Cells.Item(160, 6).Value = 100112022
$ws.Cells.Item(160, 7).Value = "Arveja Verde"
$ws.Cells.Item(160, 8).Value = "Perfection"
$ws.Cells.Item(160, 9).Value = "Primera"
$ws.Cells.Item(160, 10).Value = 52
$ws.Cells.Item(160, 11).Value = 26000
$ws.Cells.Item(160, 12).Value = 28000
$ws.Cells.Item(160, 13).Value = 27000
$ws.Cells.Item(160, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(160, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(160, 16).Value = 1080
$ws.Cells.Item(160, 17).Value = 25
$ws.Cells.Item(160, 18).Value = "Hortaliza"
